$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [string][char]0x2083
$ws.Cells.Item(2, 4).Value = '63.186.34'
$ws.Cells.Item(2, 5).Value = '  -3.02%  '
$ws.Cells.Item(3, 4).Value = '3.248.67'
$ws.Cells.Item(3, 5).Value = '  -3.98%  '
$ws.Cells.Item(4, 5).Value = '  -0.23%  '
$ws.Cells.Item(5, 4).Value = "'175.38"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -3.03%  '
$ws.Cells.Item(6, 4).Value = "'521.33"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.89%  '
$ws.Cells.Item(7, 5).Value = '  -2.24%  '
$ws.Cells.Item(8, 2).Value = 'USDC'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(8, 4).Value = "'1.00"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.02%  '
$ws.Cells.Item(9, 2).Value = 'LidoStakedEther'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Cells.Item(9, 4).Value = '3.246.33'
$ws.Cells.Item(9, 5).Value = '  -3.54%  '
$ws.Cells.Item(10, 4).Value = "'0.602"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -3.03%  '
$ws.Cells.Item(11, 4).Value = "'52.91"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -7.34%  '
$ws.Cells.Item(12, 5).Value = '  -1.30%  '
$ws.Cells.Item(13, 5).Value = '  -0.47%  '
$ws.Cells.Item(14, 4).Value = "'8.95"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -3.40%  '
$ws.Cells.Item(15, 4).Value = '3.777.37'
$ws.Cells.Item(15, 5).Value = '  -4.23%  '
$ws.Cells.Item(16, 4).Value = '3.255.61'
$ws.Cells.Item(16, 5).Value = '  -4.31%  '
$ws.Cells.Item(17, 4).Value = "'0.115"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -6.28%  '
$ws.Cells.Item(18, 4).Value = '63.168.55'
$ws.Cells.Item(18, 5).Value = '  -2.96%  '
$ws.Cells.Item(19, 4).Value = "'17.20"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -1.32%  '
$ws.Cells.Item(20, 4).Value = "'11.01"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -1.25%  '
$ws.Cells.Item(21, 4).Value = "'0.962"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.94%  '
$ws.Cells.Item(22, 4).Value = "'366.24"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -1.81%  '
$ws.Cells.Item(23, 4).Value = "'3.75"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.90%  '
$ws.Cells.Item(24, 4).Value = "'80.45"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -2.27%  '
$ws.Cells.Item(25, 4).Value = "'11.02"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +2.61%  '
$ws.Cells.Item(26, 4).Value = "'3.90"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +8.79%  '
$ws.Cells.Item(27, 5).Value = '  +4.96%  '
$ws.Cells.Item(28, 4).Value = "'2.62"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -0.52%  '
$ws.Cells.Item(29, 4).Value = "'11.25"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -2.13%  '
$ws.Cells.Item(30, 5).Value = '  -3.02%  '
$ws.Cells.Item(31, 4).Value = "'657.27"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -2.55%  '
$ws.Cells.Item(32, 4).Value = "'28.35"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -4.03%  '
$ws.Cells.Item(33, 4).Value = "'6.40"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -4.04%  '
$ws.Cells.Item(34, 4).Value = "'11.12"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.29%  '
$ws.Cells.Item(35, 4).Value = "'0.105"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -0.43%  '
$ws.Cells.Item(36, 5).Value = '  -6.24%  '
$ws.Cells.Item(37, 5).Value = '  +0.13%  '
$ws.Cells.Item(38, 4).Value = "'36.49"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +0.15%  '
$ws.Cells.Item(39, 5).Value = '  -1.86%  '
$ws.Cells.Item(40, 4).Value = "'1.00"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -0.12%  '
$ws.Cells.Item(41, 4).Value = '0.0' + $sub3 + '0715'
$ws.Cells.Item(41, 5).Value = '  +15.14%  '
$ws.Cells.Item(42, 5).Value = '  -3.37%  '
$ws.Cells.Item(43, 4).Value = '2.884.03'
$ws.Cells.Item(43, 5).Value = '  +0.75%  '
$ws.Cells.Item(44, 4).Value = "'2.49"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +7.54%  '
$ws.Cells.Item(45, 4).Value = "'2.64"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -0.35%  '
$ws.Cells.Item(46, 4).Value = "'0.0390"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.65%  '
$ws.Cells.Item(47, 5).Value = '  +10.82%  '
$ws.Cells.Item(48, 5).Value = '  -6.54%  '
$ws.Cells.Item(49, 4).Value = "'3.01"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +6.81%  '
$ws.Cells.Item(50, 4).Value = "'137.73"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +2.00%  '
$ws.Cells.Item(51, 5).Value = '  -1.59%  '
